$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = "oioi"
$ws.Range("B53").Value = "oioi"
$ws.Range("C53").Value = "oioi"
$ws.Range("D53").Value = "oioi"
$ws.Range("E53").Value = "adm"
